# Insert a new weekly price record at row 116 ("Ají" / Americana (o) from
# "Provincia de Limarí", dated 2021-09-22 / serial 44461). This pushes the
# previously existing rows 116-123 down to 117-124, extending the used
# range from A1:R123 to A1:R124, matching the committed diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 116, shifting rows 116..123 down to 117..124.
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new record's data.
$ws.Cells.Item(116, 1).Value  = 5
$ws.Cells.Item(116, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(116, 3).Value  = "Maule"
$ws.Cells.Item(116, 4).Value  = 44461
$ws.Cells.Item(116, 5).Value  = 7
$ws.Cells.Item(116, 6).Value  = 100112021
$ws.Cells.Item(116, 7).Value  = "Ají"
$ws.Cells.Item(116, 8).Value  = "Americana (o)"
$ws.Cells.Item(116, 9).Value  = "Primera"
$ws.Cells.Item(116, 10).Value = 100
$ws.Cells.Item(116, 11).Value = 75000
$ws.Cells.Item(116, 12).Value = 75000
$ws.Cells.Item(116, 13).Value = 75000
$ws.Cells.Item(116, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(116, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(116, 16).Value = 3000
$ws.Cells.Item(116, 17).Value = 25
$ws.Cells.Item(116, 18).Value = "Hortaliza"
